$wb = $excel.ActiveWorkbook

# The handback file bbad4a33-0532-423b-9acd-dcfbe91a200f.md has now come back
# in sync with en-US, so the generated report is refreshed to reflect the
# new status instead of "Ready for handoff", clearing the stale error detail
# and bumping the handback timestamps.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E6").Value = "Handed back: in sync with en-US"
$overview.Range("F6").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C6").Value = "Handed back: in sync with en-US"
$zhcn.Range("K6").Value = "2016-11-03 19:35:34"
$zhcn.Range("P6").Value = ""

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C6").Value = "Handed back: in sync with en-US"
$dede.Range("K6").Value = "2016-11-03 19:35:51"
$dede.Range("P6").Value = ""
